# Atualizacao automatica: 2025-08-29 09:00:45
# Refresh of the dashboard_data sheet: re-orders/updates the detection
# records in rows 7-11 and rows 21-22 to reflect the latest scrape.
#
# Note: columns I (First_Coords) and J (First_Confidence) are stored as
# text even though they look numeric (comma separated coordinates /
# decimal confidence scores). A leading apostrophe is used so Excel
# keeps them as text instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 7-11 (columns A, D, E, F, G, H, I, J) ---
# Row 7
$ws.Range("A7").Value = "2117575c-4ae1-458c-b88a-fc40f40debdb"
$ws.Range("D7").Value = "image_20250727074723_ppp0.jpg"
$ws.Range("E7").Value = "PLACA_20250723145134"
$ws.Range("F7").Value = "Moura"
$ws.Range("G7").Value = 38.06587
$ws.Range("H7").Value = -7.221796
$ws.Range("I7").Value = "'1490,161,1563,258"
$ws.Range("J7").Value = "'0.62"

# Row 8
$ws.Range("A8").Value = "283b6eda-9c83-4cdd-9524-c7c394f2dc89"
$ws.Range("D8").Value = "image_20250728214139_ppp0.jpg"
$ws.Range("E8").Value = "PLACA_20250717165933"
$ws.Range("F8").Value = "Beja"
$ws.Range("G8").Value = 38.02035
$ws.Range("H8").Value = -7.94715
$ws.Range("I8").Value = "'962,713,1006,765"
$ws.Range("J8").Value = "'0.76"

# Row 9
$ws.Range("A9").Value = "a19b65d1-6f97-4841-9e1c-7446a9be92b6"
$ws.Range("D9").Value = "image_20250728214139_ppp0.jpg"
$ws.Range("E9").Value = "PLACA_20250717165933"
$ws.Range("F9").Value = "Beja"
$ws.Range("G9").Value = 38.02035
$ws.Range("H9").Value = -7.94715
$ws.Range("I9").Value = "'967,614,1002,659"
$ws.Range("J9").Value = "'0.73"

# Row 10
$ws.Range("A10").Value = "4be1b1cf-d480-453e-b5fb-d4ecd6764c4d"
$ws.Range("D10").Value = "image_20250728214139_ppp0.jpg"
$ws.Range("E10").Value = "PLACA_20250717165933"
$ws.Range("F10").Value = "Beja"
$ws.Range("G10").Value = 38.02035
$ws.Range("H10").Value = -7.94715
$ws.Range("I10").Value = "'702,633,740,690"
$ws.Range("J10").Value = "'0.72"

# Row 11
$ws.Range("A11").Value = "dfd476d4-7689-4671-a076-78fe3ce806bb"
$ws.Range("D11").Value = "image_20250728214139_ppp0.jpg"
$ws.Range("E11").Value = "PLACA_20250717165933"
$ws.Range("F11").Value = "Beja"
$ws.Range("G11").Value = 38.02035
$ws.Range("H11").Value = -7.94715
$ws.Range("I11").Value = "'1254,850,1294,895"
$ws.Range("J11").Value = "'0.67"

# --- Rows 21-22 (columns A, D, I, J) ---
# Row 21
$ws.Range("A21").Value = "a2ea21b8-7dce-4e6a-be35-4edaddca5896"
$ws.Range("D21").Value = "image_20250824092407_ppp0.jpg"
$ws.Range("I21").Value = "'1002,789,1039,825"
$ws.Range("J21").Value = "'0.64"

# Row 22
$ws.Range("A22").Value = "66efa766-1456-4beb-b92a-0615a2fc41bb"
$ws.Range("D22").Value = "image_20250824214658_ppp0.jpg"
$ws.Range("I22").Value = "'1272,293,1315,331"
$ws.Range("J22").Value = "'0.69"
